$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles from column P into new column Q for the rows that have
# a formatted P cell, then set the Q values for the data rows.

# Row 3: bottom border cell only (no value)
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4: year header 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q4").Value = 2020

# Row 5: data value
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q5").Value = 25.6

# Row 6: data value
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q6").Value = 13.073527219449954

# Row 7: data value
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q7").Value = 21.941290626870046

# Row 8: total value
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q8").Value = 196.6

# Row 1 header gains a custom height
$ws.Rows("1").RowHeight = 19.5
